# Adds the "metadata" sheet (panel query metadata) after "data", and
# refreshes the "data" sheet's per-row query timestamps (column F) to the
# values captured by the re-run of the PanelApp scrape.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# ---------------------------------------------------------------------
# 1. Refresh F2:F100 ("time_taken") on the "data" sheet with the new
#    per-row query timestamps.
# ---------------------------------------------------------------------
$newTimes = @(
    "2021-10-05 14:19:23.054694",
    "2021-10-05 14:19:23.054701",
    "2021-10-05 14:19:23.054705",
    "2021-10-05 14:19:23.054707",
    "2021-10-05 14:19:23.054710",
    "2021-10-05 14:19:23.054713",
    "2021-10-05 14:19:23.054716",
    "2021-10-05 14:19:23.054718",
    "2021-10-05 14:19:23.054721",
    "2021-10-05 14:19:23.054723",
    "2021-10-05 14:19:23.054726",
    "2021-10-05 14:19:23.054729",
    "2021-10-05 14:19:23.054731",
    "2021-10-05 14:19:23.054734",
    "2021-10-05 14:19:23.054736",
    "2021-10-05 14:19:23.054739",
    "2021-10-05 14:19:23.054741",
    "2021-10-05 14:19:23.054744",
    "2021-10-05 14:19:23.054746",
    "2021-10-05 14:19:23.054749",
    "2021-10-05 14:19:23.054751",
    "2021-10-05 14:19:23.054754",
    "2021-10-05 14:19:23.054756",
    "2021-10-05 14:19:23.054759",
    "2021-10-05 14:19:23.054762",
    "2021-10-05 14:19:23.054764",
    "2021-10-05 14:19:23.054767",
    "2021-10-05 14:19:23.054769",
    "2021-10-05 14:19:23.054771",
    "2021-10-05 14:19:23.054774",
    "2021-10-05 14:19:23.054776",
    "2021-10-05 14:19:23.054779",
    "2021-10-05 14:19:23.054782",
    "2021-10-05 14:19:23.054784",
    "2021-10-05 14:19:23.054787",
    "2021-10-05 14:19:23.054789",
    "2021-10-05 14:19:23.054792",
    "2021-10-05 14:19:23.054794",
    "2021-10-05 14:19:23.054797",
    "2021-10-05 14:19:23.054799",
    "2021-10-05 14:19:23.054802",
    "2021-10-05 14:19:23.054805",
    "2021-10-05 14:19:23.054808",
    "2021-10-05 14:19:23.054810",
    "2021-10-05 14:19:23.054813",
    "2021-10-05 14:19:23.054815",
    "2021-10-05 14:19:23.054818",
    "2021-10-05 14:19:23.054820",
    "2021-10-05 14:19:23.054823",
    "2021-10-05 14:19:23.054825",
    "2021-10-05 14:19:23.054828",
    "2021-10-05 14:19:23.054830",
    "2021-10-05 14:19:23.054833",
    "2021-10-05 14:19:23.054835",
    "2021-10-05 14:19:23.054838",
    "2021-10-05 14:19:23.054840",
    "2021-10-05 14:19:23.054843",
    "2021-10-05 14:19:23.054845",
    "2021-10-05 14:19:23.054848",
    "2021-10-05 14:19:23.054850",
    "2021-10-05 14:19:23.054853",
    "2021-10-05 14:19:23.054856",
    "2021-10-05 14:19:23.054858",
    "2021-10-05 14:19:23.054861",
    "2021-10-05 14:19:23.054865",
    "2021-10-05 14:19:23.054868",
    "2021-10-05 14:19:23.054870",
    "2021-10-05 14:19:23.054873",
    "2021-10-05 14:19:23.054875",
    "2021-10-05 14:19:23.054878",
    "2021-10-05 14:19:23.054880",
    "2021-10-05 14:19:23.054883",
    "2021-10-05 14:19:23.054885",
    "2021-10-05 14:19:23.054888",
    "2021-10-05 14:19:23.054891",
    "2021-10-05 14:19:23.054893",
    "2021-10-05 14:19:23.054898",
    "2021-10-05 14:19:23.054901",
    "2021-10-05 14:19:23.054903",
    "2021-10-05 14:19:23.054906",
    "2021-10-05 14:19:23.054908",
    "2021-10-05 14:19:23.054911",
    "2021-10-05 14:19:23.054914",
    "2021-10-05 14:19:23.054916",
    "2021-10-05 14:19:23.054919",
    "2021-10-05 14:19:23.054921",
    "2021-10-05 14:19:23.054924",
    "2021-10-05 14:19:23.054926",
    "2021-10-05 14:19:23.054929",
    "2021-10-05 14:19:23.054932",
    "2021-10-05 14:19:23.054934",
    "2021-10-05 14:19:23.054937",
    "2021-10-05 14:19:23.054941",
    "2021-10-05 14:19:23.054944",
    "2021-10-05 14:19:23.054946",
    "2021-10-05 14:19:23.054949",
    "2021-10-05 14:19:23.054952",
    "2021-10-05 14:19:23.054954",
    "2021-10-05 14:19:23.054957"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $newTimes[$i]
}

# ---------------------------------------------------------------------
# 2. Add the new "metadata" sheet right after "data".
# ---------------------------------------------------------------------
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "CAKUT"
$metaSheet.Range("C2").Value = 234
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "1.164"
$metaSheet.Range("E2").Value = "2021-09-06T10:17:22.215165Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:19:23.050993"
$metaSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/234/?format=json"
